$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 173, shifting existing rows 173:205 down to 174:206
$ws.Rows.Item(173).Insert()

$ws.Cells.Item(173, 1).Value = 7
$ws.Cells.Item(173, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(173, 3).Value = "Ñuble"
$ws.Cells.Item(173, 4).Value = 44637
$ws.Cells.Item(173, 5).Value = 16
$ws.Cells.Item(173, 6).Value = "Fruta"
$ws.Cells.Item(173, 7).Value = 100108
$ws.Cells.Item(173, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(173, 9).Value = 100108005
$ws.Cells.Item(173, 10).Value = "Piña"
$ws.Cells.Item(173, 11).Value = "Caramelo"
$ws.Cells.Item(173, 12).Value = "Segunda"
$ws.Cells.Item(173, 13).Value = 60
$ws.Cells.Item(173, 14).Value = 17000
$ws.Cells.Item(173, 15).Value = 18000
$ws.Cells.Item(173, 16).Value = 17500
$ws.Cells.Item(173, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(173, 18).Value = "Ecuador"
$ws.Cells.Item(173, 19).Value = 1250
$ws.Cells.Item(173, 20).Value = 14
